$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

# Copy formatting (number format/style) from the cell above so the new
# date cell reuses the existing style index instead of minting a new one.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 42625.883391203701
$ws.Cells.Item($row, 2).Value = 32
$ws.Cells.Item($row, 3).Value = 59
$ws.Cells.Item($row, 4).Value = 36
$ws.Cells.Item($row, 5).Value = 97
$ws.Cells.Item($row, 6).Value = 2
$ws.Cells.Item($row, 7).Value = 17195
$ws.Cells.Item($row, 8).Value = 19737
$ws.Cells.Item($row, 9).Value = 2138
$ws.Cells.Item($row, 10).Value = 374
$ws.Cells.Item($row, 11).Value = 229
$ws.Cells.Item($row, 12).Value = 37
$ws.Cells.Item($row, 13).Value = 1
$ws.Cells.Item($row, 14).Value = "Bag"
